# corrected original_id in metadata
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = "o2885712-ME7_hyp"
$ws.Range("A12").Value = "o2885713-mNS_hyp"
$ws.Range("A13").Value = "o2885714-RML6_hyp"
